$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "TEAM TABLE (HR)" field list (H16:I20) gets a new row inserted at the
# top ("Agent ID" / "int"), pushing the existing rows down by one and
# correcting a few data types along the way. The old note in J18 is removed.

# Remove the old "*is there a phone data type?" note next to Office Phone.
$ws.Range("J18").ClearContents()

# Shift existing rows down (bottom-up, so we don't clobber data before it
# is copied).
$ws.Range("H21").Value = $ws.Range("H20").Value2
$ws.Range("I21").Value = "text"

$ws.Range("H20").Value = $ws.Range("H19").Value2
$ws.Range("I20").Value = "text"

$ws.Range("H19").Value = $ws.Range("H18").Value2
$ws.Range("I19").Value = "varchar"

$ws.Range("H18").Value = $ws.Range("H17").Value2
$ws.Range("I18").Value = "text"

$ws.Range("H17").Value = $ws.Range("H16").Value2
$ws.Range("I17").Value = "text"

# New first row of the field list.
$ws.Range("H16").Value = "Agent ID"
$ws.Range("I16").Value = "int"

# Highlight the newly added row with a white fill.
$ws.Range("H16").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1

# Restore the selection to what it was left at in the edited workbook.
$ws.Range("J25").Select()
